$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need a leading apostrophe
# to stop Excel from auto-converting them to a numeric value (which would
# lose the original string formatting / exact digits, e.g. trailing zeros).
# The Style reset afterwards clears the "quote prefix" formatting flag that
# the apostrophe entry implicitly applies, so the cell keeps its original
# (unstyled) look exactly like the rest of the sheet.

$ws.Range("D2").Value = '64.075.50'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.738.92'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'570.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").Value = "'160.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = "'0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.73%  '
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("E10").Value = '  +4.63%  '
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").Value = '3.222.31'
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").Value = "'26.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '63.891.27'
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").Value = '2.746.40'
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").Value = "'12.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").Value = "'4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").Value = "'354.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E23").Value = '  -5.25%  '
$ws.Range("D24").Value = "'64.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = "'8.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("D28").Value = '0.0₃0915'
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = "'7.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.94%  '
$ws.Range("E31").Value = '  +8.05%  '
$ws.Range("D32").Value = "'163.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.11%  '
$ws.Range("D33").Value = "'4.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").Value = "'20.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.40%  '
$ws.Range("E35").Value = '  +1.89%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("D38").Value = "'0.993"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").Value = "'351.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.16%  '
$ws.Range("D40").Value = "'6.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.01%  '
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D43").Value = "'22.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.40%  '
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("D45").Value = "'0.0585"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.50%  '
$ws.Range("E46").Value = '  -1.42%  '
$ws.Range("D47").Value = "'134.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.30%  '
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").Value = "'0.0250"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.48%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = "'11.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.10%  '
